$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte Consolidado")

# Header timestamp update
$ws.Range("A4").Value = "Generado: 13/07/2025 08:06 | Total de Tickets: 146"

# Row 14/15: swap K column values
$ws.Range("K14").Value = "resuelto"
$ws.Range("K15").Value = "ticket resuelto"

# Row 40: ticket resolved with notes
$ws.Range("D40").Value = "cerrado"
$ws.Range("G40").Value = "12/07/2025 23:18"
$ws.Range("J40").Value = "12/07/2025 23:18"
$ws.Range("K40").Value = "TEST RESOLUTION NOTES for TKT-000106 - This should appear in email notifications"

# Row 41/42: swap K column values
$ws.Range("K41").Value = "resolucion 2"
$ws.Range("K42").Value = "ticket resuelto"

# Row 76: ticket closed
$ws.Range("D76").Value = "cerrado"
$ws.Range("G76").Value = "12/07/2025 23:19"
$ws.Range("J76").Value = "12/07/2025 23:19"
$ws.Range("K76").Value = "prueba de que si se resolvvio ticket 65"

# Row 77: ticket closed
$ws.Range("D77").Value = "cerrado"
$ws.Range("G77").Value = "12/07/2025 23:27"
$ws.Range("J77").Value = "12/07/2025 23:27"
$ws.Range("K77").Value = "resolucion ticket 64"

# Rows 144-147: shift K column values down by one (cyclic-ish)
# K144/K145 are numeric-looking text ("34343434" / "1") - prefix with an
# apostrophe so Excel stores them as text (matching the source data),
# not auto-converted to numbers.
$ws.Range("K144").Value = "'34343434"
$ws.Range("K145").Value = "'1"
$ws.Range("K146").Value = "prueba de resolucion 40"
$ws.Range("K147").Value = "sdsds"
